# Generate Report for Handback
#
# This localization-status report records, per language, the latest
# handoff/handback state for each source file. This run represents a
# handback event: the overall status moves from "Ready for handoff" to
# "Handed back: in sync with en-US", the handback timestamp is recorded,
# and the "Latest Target File" / "Latest Handback File" columns (which
# were previously empty placeholders) are now populated with links to the
# generated target/handback artifact, mirroring the existing handoff-file
# hyperlinks already present on each row.

$wb = $excel.ActiveWorkbook

# --- 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# The Overview sheet's B/C columns mirror each language sheet's Status (C)
# column, so update all of them together.
$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# --- 2. Populate "Latest Target File" (F) / "Latest Handback File" (G) ---
# These link out to the same target/handback artifacts already referenced
# by column D (Latest Handoff File) on each row.
$zhcnTargetUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/19bc1a1d6b0d29df505ad7ff09f4c48a1ee18ffa/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcnTargetName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$aMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/a8cead93de34f3dbd422630230b70566dee5d740/e2e/a.md"

$zhcn.Hyperlinks.Add($zhcn.Range("F2"), $aMdUrl, "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), $zhcnTargetUrl, "", "", $zhcnTargetName)
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), $aMdUrl, "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), $zhcnTargetUrl, "", "", $zhcnTargetName)

$dedeTargetUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e6d4adec54cd9d7e7d0c814042a2932650e7ddb5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dedeTargetName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$dede.Hyperlinks.Add($dede.Range("F2"), $aMdUrl, "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("G2"), $dedeTargetUrl, "", "", $dedeTargetName)
$dede.Hyperlinks.Add($dede.Range("F3"), $aMdUrl, "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("G3"), $dedeTargetUrl, "", "", $dedeTargetName)

# --- 3. Latest Handback DateTime (H): record the handback timestamps ---
# zh-cn finished syncing slightly before de-de.
$zhcn.Range("H2").Value = "2016-03-22 02:28:34"
$zhcn.Range("H3").Value = "2016-03-22 02:28:34"
$dede.Range("H2").Value = "2016-03-22 02:28:41"
$dede.Range("H3").Value = "2016-03-22 02:28:41"
